# Weekly cryptos-list refresh (GitHub Actions bot) -- Sat Apr 29 23:07:15 UTC 2023
# Updates the "Price" (D) and "Volume(1h)" (E) columns for rows 2-51 to match
# the latest scrape. Values are plain text (not numbers), so numeric-looking
# prices are written with a leading apostrophe to force Excel to keep them as
# text instead of silently re-parsing/rounding them (e.g. "1.010" -> 1.01).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.396.73"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.916.14"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.64%  "
$ws.Range("D5").Value = "`'324.82"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "`'1.007"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "`'0.4813"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").Value = "`'0.4059"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "`'0.08201"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").Value = "`'1.010"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("D11").Value = "`'23.18"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "1.902.55"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "`'6.062"
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").Value = "`'7.223"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").Value = "`'91.57"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").Value = "`'0.06858"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "`'0.00001038"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").Value = "`'17.57"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "`'1.007"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "29.420.50"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "`'5.658"
$ws.Range("E22").Value = "  +2.16%  "
$ws.Range("D23").Value = "`'11.72"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("D24").Value = "`'2.192"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").Value = "2.147.00"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").Value = "`'6.592"
$ws.Range("E26").Value = "  +8.09%  "
$ws.Range("D27").Value = "`'155.90"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").Value = "`'19.97"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").Value = "`'2.110"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").Value = "`'120.59"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("D31").Value = "`'1.014"
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("D32").Value = "`'0.09616"
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("D33").Value = "`'5.649"
$ws.Range("E33").Value = "  +3.23%  "
$ws.Range("D34").Value = "`'3.544"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "`'1.371"
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("D36").Value = "`'0.02283"
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("D37").Value = "`'0.06099"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("D38").Value = "`'1.181"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").Value = "`'10.88"
$ws.Range("E39").Value = "  +6.80%  "
$ws.Range("D40").Value = "`'8.059"
$ws.Range("E40").Value = "  +1.81%  "
$ws.Range("D41").Value = "`'0.5959"
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("D42").Value = "`'0.1843"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "`'1.279"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").Value = "`'2.377"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").Value = "`'0.07610"
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("D46").Value = "`'12.46"
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("D47").Value = "`'0.5581"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("D48").Value = "`'1.949"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").Value = "`'118.40"
$ws.Range("E49").Value = "  +3.96%  "
$ws.Range("D50").Value = "`'2.430"
$ws.Range("E50").Value = "  +3.85%  "
$ws.Range("D51").Value = "`'72.14"
$ws.Range("E51").Value = "  -0.31%  "
